$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 10.47960307264752
$ws.Range("D2").Value = 3.83824738459685
$ws.Range("E2").Value = 12.76766712996633
$ws.Range("F2").Value = 24.38713471282597
$ws.Range("G2").Value = 30.25043463407877
$ws.Range("H2").Value = 13.72114037120475
$ws.Range("I2").Value = 22.53310231119853
$ws.Range("L2").Value = 9.329841130967061
$ws.Range("O2").Value = 21.2833099749091
$ws.Range("C3").Value = 10.45037202107952
$ws.Range("D3").Value = 3.835930178954499
$ws.Range("E3").Value = 12.72644113586963
$ws.Range("F3").Value = 24.05226021884549
$ws.Range("G3").Value = 29.53861713261317
$ws.Range("H3").Value = 13.68306165960474
$ws.Range("I3").Value = 22.35737378144347
$ws.Range("L3").Value = 9.328454093328642
$ws.Range("O3").Value = 21.09573507956767
$ws.Range("C4").Value = 10.43465040749228
$ws.Range("D4").Value = 3.834424825415421
$ws.Range("E4").Value = 12.70398083294631
$ws.Range("F4").Value = 23.85142729026402
$ws.Range("G4").Value = 29.10242667723949
$ws.Range("H4").Value = 13.66258951676806
$ws.Range("I4").Value = 22.25470344097577
$ws.Range("L4").Value = 9.32929501348675
$ws.Range("O4").Value = 20.98539993083244
$ws.Range("C5").Value = 10.42880791258606
$ws.Range("D5").Value = 3.833790664373361
$ws.Range("E5").Value = 12.69555162918325
$ws.Range("F5").Value = 23.77089341594906
$ws.Range("G5").Value = 28.92518720548114
$ws.Range("H5").Value = 13.65498394736394
$ws.Range("I5").Value = 22.21421899936945
$ws.Range("L5").Value = 9.330063813335169
$ws.Range("O5").Value = 20.94169908310523
$ws.Range("C6").Value = 10.42787196222097
$ws.Range("D6").Value = 3.833684112561664
$ws.Range("E6").Value = 12.69419583309895
$ws.Range("F6").Value = 23.75760278919316
$ws.Range("G6").Value = 28.89579634248949
$ws.Range("H6").Value = 13.6537657119091
$ws.Range("I6").Value = 22.20757951812481
$ws.Range("L6").Value = 9.330217219497595
$ws.Range("O6").Value = 20.93452006577358
$ws.Range("C7").Value = 10.43456932385896
$ws.Range("D7").Value = 3.83441635672833
$ws.Range("E7").Value = 12.70386421649433
$ws.Range("F7").Value = 23.85033575306452
$ws.Range("G7").Value = 29.10003389911549
$ws.Range("H7").Value = 13.66248395444741
$ws.Range("I7").Value = 22.25415191823619
$ws.Range("L7").Value = 9.32930365592838
$ws.Range("O7").Value = 20.98480539963697
$ws.Range("C8").Value = 10.46906490508412
$ws.Range("D8").Value = 3.837465579923371
$ws.Range("E8").Value = 12.75286433203479
$ws.Range("F8").Value = 24.27074128271625
$ws.Range("G8").Value = 30.00499325919384
$ws.Range("H8").Value = 13.70741082898344
$ws.Range("I8").Value = 22.47145089682965
$ws.Range("L8").Value = 9.329011947853191
$ws.Range("O8").Value = 21.21765833673252
$ws.Range("C9").Value = 10.55415950457158
$ws.Range("D9").Value = 3.842787855798174
$ws.Range("E9").Value = 12.87127721233442
$ws.Range("F9").Value = 25.12838031552726
$ws.Range("G9").Value = 31.77406878043188
$ws.Range("H9").Value = 13.81830698560665
$ws.Range("I9").Value = 22.937175097115
$ws.Range("L9").Value = 9.341840360411121
$ws.Range("O9").Value = 21.71054551788778
$ws.Range("C10").Value = 10.62699022019643
$ws.Range("D10").Value = 3.846294917704664
$ws.Range("E10").Value = 12.97141709788208
$ws.Range("F10").Value = 25.77235484085015
$ws.Range("G10").Value = 33.05458916311124
$ws.Range("H10").Value = 13.91325508847762
$ws.Range("I10").Value = 23.30081928520842
$ws.Range("L10").Value = 9.359385048239247
$ws.Range("O10").Value = 22.09180496302789
$ws.Range("C11").Value = 10.66228121176307
$ws.Range("D11").Value = 3.847802294688922
$ws.Range("E11").Value = 13.01971159451615
$ws.Range("F11").Value = 26.06697412088439
$ws.Range("G11").Value = 33.6298858845403
$ws.Range("H11").Value = 13.9592700100001
$ws.Range("I11").Value = 23.47031043443012
$ws.Range("L11").Value = 9.369114970458863
$ws.Range("O11").Value = 22.26874091314613
$ws.Range("C12").Value = 10.67594794756043
$ws.Range("D12").Value = 3.848360398001666
$ws.Range("E12").Value = 13.03838266102604
$ws.Range("F12").Value = 26.17866826080345
$ws.Range("G12").Value = 33.84646349813307
$ws.Range("H12").Value = 13.97709061067314
$ws.Range("I12").Value = 23.53502473213244
$ws.Range("L12").Value = 9.373049282135547
$ws.Range("O12").Value = 22.33618803746316
$ws.Range("C13").Value = 10.67299123339563
$ws.Range("D13").Value = 3.848240766474034
$ws.Range("E13").Value = 13.03434466247179
$ws.Range("F13").Value = 26.15460898611309
$ws.Range("G13").Value = 33.79987991998833
$ws.Range("H13").Value = 13.97323519235058
$ws.Range("I13").Value = 23.52106455949447
$ws.Range("L13").Value = 9.372190876534312
$ws.Range("O13").Value = 22.3216432121414
$ws.Range("C14").Value = 10.663399560055
$ws.Range("D14").Value = 3.847848465165526
$ws.Range("E14").Value = 13.02124006446527
$ws.Range("F14").Value = 26.07616132903326
$ws.Range("G14").Value = 33.64773069960466
$ws.Range("H14").Value = 13.96072825039585
$ws.Range("I14").Value = 23.47562417072504
$ws.Range("L14").Value = 9.369433652557957
$ws.Range("O14").Value = 22.27428123662491
$ws.Range("C15").Value = 10.65756357761039
$ws.Range("D15").Value = 3.847606512817661
$ws.Range("E15").Value = 13.01326265647562
$ws.Range("F15").Value = 26.02812328403852
$ws.Range("G15").Value = 33.55436215727691
$ws.Range("H15").Value = 13.95311861412588
$ws.Range("I15").Value = 23.44785826721343
$ws.Range("L15").Value = 9.367777251676774
$ws.Range("O15").Value = 22.24532692856566
$ws.Range("C16").Value = 10.62472667925059
$ws.Range("D16").Value = 3.846194631687089
$ws.Range("E16").Value = 12.96831515068649
$ws.Range("F16").Value = 25.75312477416255
$ws.Range("G16").Value = 33.0168262926481
$ws.Range("H16").Value = 13.9103038504742
$ws.Range("I16").Value = 23.28981995863522
$ws.Range("L16").Value = 9.35878422663583
$ws.Range("O16").Value = 22.0803072598365
$ws.Range("C17").Value = 10.60513012757232
$ws.Range("D17").Value = 3.845305910785877
$ws.Range("E17").Value = 12.94143551839098
$ws.Range("F17").Value = 25.58476974013746
$ws.Range("G17").Value = 32.68504285504711
$ws.Range("H17").Value = 13.88475429673296
$ws.Range("I17").Value = 23.19387278324225
$ws.Range("L17").Value = 9.353714075845884
$ws.Range("O17").Value = 21.97992825560868
$ws.Range("C18").Value = 10.59406251079553
$ws.Range("D18").Value = 3.84478646184106
$ws.Range("E18").Value = 12.92623370685781
$ws.Range("F18").Value = 25.48810033995068
$ws.Range("G18").Value = 32.4935434718583
$ws.Range("H18").Value = 13.87032518876805
$ws.Range("I18").Value = 23.13907261654626
$ws.Range("L18").Value = 9.350962573605361
$ws.Range("O18").Value = 21.92252576625884
$ws.Range("C19").Value = 10.59035043845753
$ws.Range("D19").Value = 3.844609164859149
$ws.Range("E19").Value = 12.92113138313663
$ws.Range("F19").Value = 25.45540136201317
$ws.Range("G19").Value = 32.42859811871104
$ws.Range("H19").Value = 13.86548577975244
$ws.Range("I19").Value = 23.1205861446342
$ws.Range("L19").Value = 9.350059299437795
$ws.Range("O19").Value = 21.90314919282136
$ws.Range("C20").Value = 10.60719517618118
$ws.Range("D20").Value = 3.84540137398612
$ws.Range("E20").Value = 12.94427021076845
$ws.Range("F20").Value = 25.60267527220191
$ws.Range("G20").Value = 32.72043249904236
$ws.Range("H20").Value = 13.88744659862062
$ws.Range("I20").Value = 23.20404694578683
$ws.Range("L20").Value = 9.354236766142611
$ws.Range("O20").Value = 21.99057973980204
$ws.Range("C21").Value = 10.66620871351043
$ws.Range("D21").Value = 3.847964038868485
$ws.Range("E21").Value = 13.02507889982562
$ws.Range("F21").Value = 26.09920070939484
$ws.Range("G21").Value = 33.69245701387616
$ws.Range("H21").Value = 13.96439118574072
$ws.Range("I21").Value = 23.4889571129236
$ws.Range("L21").Value = 9.370236750873936
$ws.Range("O21").Value = 22.28818097415049
$ws.Range("C22").Value = 10.70653849415418
$ws.Range("D22").Value = 3.849564794649631
$ws.Range("E22").Value = 13.08011881599116
$ws.Range("F22").Value = 26.42440601491083
$ws.Range("G22").Value = 34.32020302926982
$ws.Range("H22").Value = 14.01698051735751
$ws.Range("I22").Value = 23.67823623155889
$ws.Range("L22").Value = 9.382148777358454
$ws.Range("O22").Value = 22.48524919821001
$ws.Range("C23").Value = 10.6848552943018
$ws.Range("D23").Value = 3.848717238062493
$ws.Range("E23").Value = 13.05054305136375
$ws.Range("F23").Value = 26.25081040271164
$ws.Range("G23").Value = 33.98592584568908
$ws.Range("H23").Value = 13.98870551089903
$ws.Range("I23").Value = 23.57695091822163
$ws.Range("L23").Value = 9.375658567017284
$ws.Range("O23").Value = 22.37985435241051
$ws.Range("C24").Value = 10.60626094819072
$ws.Range("D24").Value = 3.845358241561818
$ws.Range("E24").Value = 12.94298786166992
$ws.Range("F24").Value = 25.5945797982112
$ws.Range("G24").Value = 32.70443516708632
$ws.Range("H24").Value = 13.88622859925332
$ws.Range("I24").Value = 23.19944607785136
$ws.Range("L24").Value = 9.353999948776787
$ws.Range("O24").Value = 21.98576324384745
$ws.Range("C25").Value = 10.52930253894838
$ws.Range("D25").Value = 3.841418667015757
$ws.Range("E25").Value = 12.83689861319032
$ws.Range("F25").Value = 24.89346810154541
$ws.Range("G25").Value = 31.29776813676036
$ws.Range("H25").Value = 13.78590716680931
$ws.Range("I25").Value = 22.80721247484069
$ws.Range("L25").Value = 9.336939626364597
$ws.Range("O25").Value = 21.57362437670841
